$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P1").Value = "Thời gian bắt đầu gia hạn giờ ra"
$ws.Range("Q1").Value = "Thời gian kết thúc gia hạn giờ ra"
$ws.Range("R1").Value = "Ngưỡng trễ sớm giờ ra"

$ws.Range("R2").Select()
